$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.99770000000001
$ws.Range("C3").Value = -12.10739999999999
$ws.Range("D19").Value = -8.523199999999997
$ws.Range("A21").Value = -20.26209999999999
$ws.Range("A23").Value = -20.12899999999998
$ws.Range("C24").Value = -13.31549999999999
$ws.Range("D24").Value = -7.890099999999999
$ws.Range("A25").Value = -21.58909999999999
$ws.Range("B27").Value = 5.983100000000007
$ws.Range("D30").Value = -7.221100000000003
$ws.Range("B31").Value = 5.627100000000002
$ws.Range("D31").Value = -8.136799999999999
$ws.Range("D33").Value = -8.068399999999999
$ws.Range("B39").Value = 10.0248
$ws.Range("B48").Value = 5.296500000000004
$ws.Range("B51").Value = 5.426900000000002
$ws.Range("B52").Value = 5.2014
$ws.Range("A53").Value = -21.89320000000001
$ws.Range("B55").Value = 6.246999999999996
$ws.Range("D55").Value = -7.912599999999998
$ws.Range("B56").Value = 5.137099999999998
$ws.Range("A57").Value = -22.0057
$ws.Range("B57").Value = 5.425699999999998
$ws.Range("C57").Value = -13.22839999999999
$ws.Range("A59").Value = -22.4051
$ws.Range("C61").Value = -13.25439999999999
$ws.Range("D65").Value = -8.177099999999999
$ws.Range("A69").Value = -21.63799999999999
$ws.Range("C70").Value = -12.2266
$ws.Range("D70").Value = -8.191799999999999
$ws.Range("B73").Value = 8.512699999999997
$ws.Range("D75").Value = -8.236000000000002
$ws.Range("A79").Value = -20.89410000000002
$ws.Range("A83").Value = -22.13899999999999
$ws.Range("D83").Value = -8.357799999999996
$ws.Range("C86").Value = -13.17539999999999
$ws.Range("B89").Value = 5.391299999999993
$ws.Range("B90").Value = 5.626999999999999
$ws.Range("A93").Value = -21.3055
$ws.Range("D96").Value = -7.3849
$ws.Range("D97").Value = -8.369100000000001
$ws.Range("C98").Value = -12.07109999999999
$ws.Range("C100").Value = -12.22459999999999
$ws.Range("C102").Value = -13.64630000000001
